$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- 1. Change number format (General) on existing D2:D4 cells so they pick up
#        the new "numFmtId=0, border, left/top/wrap" style (becomes cellXfs index 8) ---
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D4").NumberFormat = "General"

# --- 2. Row heights ---
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 45

# --- 3. New header cells G1:I1 (copy header style from F1) ---
$ws.Range("F1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)

# --- 4. New columns G:I for rows 2 through 11 (copy centered/wrap style from B4) ---
$ws.Range("B4").Copy()
$ws.Range("G2:I11").PasteSpecial(-4122)

# --- 5. New row 5 A:F : copy formatting from row 4 (already using the updated D-style) ---
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)
$ws.Range("A5:F5").ClearContents()
$ws.Rows.Item(5).RowHeight = 30

# --- 6. Fill in the shared-string values, in the same order the original author typed them ---
$ws.Range("A5").Value = "[IDS].[get_view_outgoing_cars_of_id_sostav]"
$ws.Range("B5").Value = "int @id_sostav"
$ws.Range("D5").Value = "Получить полную информацию по вагонам отправляемого состава"
$ws.Range("G5").Value = "api/ids/rwt/outgoing_cars/view/sostav/id/"
$ws.Range("H5").Value = "IDS_RWT_OutgoingCarsController"
$ws.Range("I5").Value = "ids_wsd.prototype.getViewOutgoingCarsOfIDSostav"
$ws.Range("G1").Value = "api"
$ws.Range("H1").Value = "Controller"
$ws.Range("I1").Value = "js module"
$ws.Range("F5").Value = "Форма АРМ (отправить на УЗ), Отчеты по отправке."
$ws.Range("C5").Value = "таблица"

# --- 7. Empty rows 6:11, columns A:F (copy formatting from row 5, keep D reverted to old style) ---
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F11").PasteSpecial(-4122)
$ws.Range("A6:F11").ClearContents()
$ws.Range("F2").Copy()
$ws.Range("D6:D11").PasteSpecial(-4122)

# --- 8. Column widths ---
$ws.Columns.Item(1).ColumnWidth = 50.4518229166667
$ws.Columns.Item(4).ColumnWidth = 52.4518229166667
$ws.Columns.Item(7).ColumnWidth = 39.3072916666667
$ws.Columns.Item(8).ColumnWidth = 30.5924479166667
$ws.Columns.Item(9).ColumnWidth = 47.5924479166667

# --- 9. Selection / view ---
$ws.Activate()
$ws.Range("F7").Select()

# --- 10. Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
